$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 1.75
$ws.Range("H2").Value = 3.2
$ws.Range("I2").Value = 5.5
$ws.Range("J2").Value = 2.6
$ws.Range("K2").Value = 1.91
$ws.Range("L2").Value = 6.5
$ws.Range("M2").Value = 1.14
$ws.Range("N2").Value = 5.5
$ws.Range("O2").Value = 1.62
$ws.Range("P2").Value = 2.2
$ws.Range("S2").Value = 1.62
$ws.Range("T2").Value = 2.2
$ws.Range("W2").Value = 4.5
$ws.Range("X2").Value = 6.5
$ws.Range("Z2").Value = 13
$ws.Range("AD2").Value = 7
$ws.Range("AE2").Value = 26
$ws.Range("AF2").Value = 126
$ws.Range("AG2").Value = 9.5
$ws.Range("AI2").Value = 21
$ws.Range("AJ2").Value = 67
$ws.Range("AN2").Value = 3.5
$ws.Range("AO2").Value = 10
$ws.Range("AT2").Value = 2.2
$ws.Range("AW2").Value = 7
$ws.Range("AY2").Value = 51
$ws.Range("AZ2").Value = 151
$ws.Range("BA2").Value = 251

# Row 4 updates
$ws.Range("G4").Value = 1.73
$ws.Range("H4").Value = 3.1
$ws.Range("M4").Value = 1.14
$ws.Range("N4").Value = 5.5
$ws.Range("O4").Value = 1.62
$ws.Range("P4").Value = 2.2
$ws.Range("S4").Value = 1.67
$ws.Range("T4").Value = 2.1
$ws.Range("Y4").Value = 8
$ws.Range("AO4").Value = 10
$ws.Range("AQ4").Value = 41
$ws.Range("AT4").Value = 2.1
